# Auto-generated edit script applying numeric updates to the Leve profit tables
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 62503484
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 62503484
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 187510452
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -187512200
$ws.Range("H72").Value = 62503484
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 62503484
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 562531356
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -562540092
$ws.Range("H112").Value = 4643.3335
$ws.Range("J112").Value = 4871.3726
$ws.Range("L112").Value = 14614.1178
$ws.Range("N112").Value = -16830.1178
$ws.Range("H132").Value = 2088.8909
$ws.Range("I132").Value = 2026.4117
$ws.Range("J132").Value = 2885.5
$ws.Range("K132").Value = 6079.2351
$ws.Range("L132").Value = 8656.5
$ws.Range("M132").Value = -3549.2351
$ws.Range("N132").Value = -13716.5
$ws.Range("H138").Value = 4378.791
$ws.Range("I138").Value = 4504.857
$ws.Range("J138").Value = 4355.87
$ws.Range("K138").Value = 13514.571
$ws.Range("L138").Value = 13067.61
$ws.Range("M138").Value = -8374.571
$ws.Range("N138").Value = -23347.61

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 25018.5
$ws.Range("J37").Value = 25018.5
$ws.Range("L37").Value = 25018.5
$ws.Range("N37").Value = -25564.5
$ws.Range("H44").Value = 19999
$ws.Range("J44").Value = 19999
$ws.Range("L44").Value = 19999
$ws.Range("N44").Value = -20975
$ws.Range("H45").Value = 1720.1786
$ws.Range("I45").Value = 1125.8334
$ws.Range("J45").Value = 2790
$ws.Range("K45").Value = 1125.8334
$ws.Range("L45").Value = 2790
$ws.Range("M45").Value = -748.8334
$ws.Range("N45").Value = -3544
$ws.Range("H46").Value = 5002479
$ws.Range("J46").Value = 5558228
$ws.Range("L46").Value = 5558228
$ws.Range("N46").Value = -5558866
$ws.Range("H55").Value = 19999
$ws.Range("J55").Value = 19999
$ws.Range("L55").Value = 19999
$ws.Range("N55").Value = -20629
$ws.Range("H61").Value = 7409988
$ws.Range("I61").Value = 10102654
$ws.Range("J61").Value = 5157.5
$ws.Range("K61").Value = 10102654
$ws.Range("L61").Value = 5157.5
$ws.Range("M61").Value = -10102442
$ws.Range("N61").Value = -5581.5
$ws.Range("H80").Value = 20570.572
$ws.Range("J80").Value = 20570.572
$ws.Range("L80").Value = 20570.572
$ws.Range("N80").Value = -22566.572
$ws.Range("H83").Value = 20570.572
$ws.Range("J83").Value = 20570.572
$ws.Range("L83").Value = 61711.716
$ws.Range("N83").Value = -71695.716
$ws.Range("H122").Value = 46365.727
$ws.Range("I122").Value = 63398.75
$ws.Range("J122").Value = 944.3333
$ws.Range("K122").Value = 190196.25
$ws.Range("L122").Value = 2832.9999
$ws.Range("M122").Value = -187746.25
$ws.Range("N122").Value = -7732.9999
$ws.Range("H136").Value = 7409988
$ws.Range("I136").Value = 10102654
$ws.Range("J136").Value = 5157.5
$ws.Range("K136").Value = 30307962
$ws.Range("L136").Value = 15472.5
$ws.Range("M136").Value = -30305412
$ws.Range("N136").Value = -20572.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 38463390
$ws.Range("I20").Value = 1789.8889
$ws.Range("J20").Value = 125001990
$ws.Range("K20").Value = 1789.8889
$ws.Range("L20").Value = 125001990
$ws.Range("M20").Value = -1542.8889
$ws.Range("N20").Value = -125002484
$ws.Range("H75").Value = 20116.053
$ws.Range("I75").Value = 3897.8
$ws.Range("J75").Value = 25908.285
$ws.Range("K75").Value = 3897.8
$ws.Range("L75").Value = 25908.285
$ws.Range("M75").Value = -2961.8
$ws.Range("N75").Value = -27780.285
$ws.Range("H78").Value = 20116.053
$ws.Range("I78").Value = 3897.8
$ws.Range("J78").Value = 25908.285
$ws.Range("K78").Value = 11693.4
$ws.Range("L78").Value = 77724.855
$ws.Range("M78").Value = -7013.400000000001
$ws.Range("N78").Value = -87084.855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7647.494
$ws.Range("I31").Value = 1215.4375
$ws.Range("J31").Value = 9138.985000000001
$ws.Range("K31").Value = 1215.4375
$ws.Range("L31").Value = 9138.985000000001
$ws.Range("M31").Value = -920.4375
$ws.Range("N31").Value = -9728.985000000001
$ws.Range("H34").Value = 7647.494
$ws.Range("I34").Value = 1215.4375
$ws.Range("J34").Value = 9138.985000000001
$ws.Range("K34").Value = 1215.4375
$ws.Range("L34").Value = 9138.985000000001
$ws.Range("M34").Value = -1013.4375
$ws.Range("N34").Value = -9542.985000000001
$ws.Range("H58").Value = 2967.6667
$ws.Range("I58").Value = 2701.5
$ws.Range("K58").Value = 2701.5
$ws.Range("M58").Value = -2498.5
$ws.Range("H99").Value = 2404.0908
$ws.Range("I99").Value = 2176
$ws.Range("J99").Value = 2518.1365
$ws.Range("K99").Value = 2176
$ws.Range("L99").Value = 2518.1365
$ws.Range("M99").Value = -678
$ws.Range("N99").Value = -5514.136500000001
$ws.Range("H122").Value = 2150.5
$ws.Range("I122").Value = 1296.4
$ws.Range("J122").Value = 2435.2
$ws.Range("K122").Value = 3889.2
$ws.Range("L122").Value = 7305.599999999999
$ws.Range("M122").Value = -1439.2
$ws.Range("N122").Value = -12205.6
$ws.Range("H126").Value = 2404.0908
$ws.Range("I126").Value = 2176
$ws.Range("J126").Value = 2518.1365
$ws.Range("K126").Value = 6528
$ws.Range("L126").Value = 7554.4095
$ws.Range("M126").Value = -4058
$ws.Range("N126").Value = -12494.4095
$ws.Range("H136").Value = 2967.6667
$ws.Range("I136").Value = 2701.5
$ws.Range("K136").Value = 8104.5
$ws.Range("M136").Value = -5554.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1827.2
$ws.Range("I55").Value = 1002
$ws.Range("J55").Value = 1954.1538
$ws.Range("K55").Value = 3006
$ws.Range("L55").Value = 5862.4614
$ws.Range("M55").Value = -2829
$ws.Range("N55").Value = -6216.4614
$ws.Range("H113").Value = 617.35297
$ws.Range("I113").Value = 590.36957
$ws.Range("K113").Value = 1771.10871
$ws.Range("M113").Value = 398.89129
$ws.Range("H117").Value = 2850.5454
$ws.Range("I117").Value = 2452
$ws.Range("J117").Value = 3000
$ws.Range("K117").Value = 7356
$ws.Range("L117").Value = 9000
$ws.Range("M117").Value = -3914
$ws.Range("N117").Value = -15884
$ws.Range("H129").Value = 1319116.4
$ws.Range("J129").Value = 1685379.5
$ws.Range("L129").Value = 5056138.5
$ws.Range("N129").Value = -5066138.5
$ws.Range("H131").Value = 3048.8032
$ws.Range("J131").Value = 3332.9075
$ws.Range("L131").Value = 9998.7225
$ws.Range("N131").Value = -20078.7225
$ws.Range("H138").Value = 2773.875
$ws.Range("I138").Value = 1587.6666
$ws.Range("J138").Value = 6332.5
$ws.Range("K138").Value = 4762.9998
$ws.Range("L138").Value = 18997.5
$ws.Range("M138").Value = 377.0002000000004
$ws.Range("N138").Value = -29277.5
$ws.Range("H140").Value = 1821.775
$ws.Range("I140").Value = 1593.125
$ws.Range("K140").Value = 4779.375
$ws.Range("M140").Value = 400.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2074.75
$ws.Range("I102").Value = 1920.3334
$ws.Range("K102").Value = 1920.3334
$ws.Range("M102").Value = -298.3334
$ws.Range("H126").Value = 3704.5454
$ws.Range("I126").Value = 3418.3635
$ws.Range("K126").Value = 10255.0905
$ws.Range("M126").Value = -7785.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 88602
$ws.Range("J103").Value = 88602
$ws.Range("L103").Value = 88602
$ws.Range("N103").Value = -90946
$ws.Range("H122").Value = 5911.0835
$ws.Range("I122").Value = 4364.2856
$ws.Range("J122").Value = 6895.409
$ws.Range("K122").Value = 13092.8568
$ws.Range("L122").Value = 20686.227
$ws.Range("M122").Value = -10642.8568
$ws.Range("N122").Value = -25586.227

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5611272.5
$ws.Range("I132").Value = 2286.9092
$ws.Range("J132").Value = 15353195
$ws.Range("K132").Value = 6860.7276
$ws.Range("L132").Value = 46059585
$ws.Range("M132").Value = -4330.7276
$ws.Range("N132").Value = -46064645
